$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (in the order values are provided per row) mapped to column index
$cols = @(5, 7, 8, 9, 10, 11, 13, 14, 15, 16, 17, 18, 19, 20)

# New values per row r=2..21, in the same column order as $cols
$newValues = @{}
$newValues[2] = @(3, 5.349936666666667, 16.04981, 0.04148245374559899, 0.04148245374559899, 3, 33.67810466666666, 101.034314, 0.1611708087274131, 0.1681783058563055, 180.1757270200378, 1621.58154318034, 0.006685760618175698, 0.006976448793697392)
$newValues[3] = @(3, 5.349936666666667, 16.04981, 0.04148245374559899, 0.04148245374559899, 3, 51.86084433333334, 155.582533, 0.2481865979460151, 0.2589774284088547, 277.4522326631923, 2497.07009396873, 0.01029538906957315, 0.01074301919512449)
$newValues[4] = @(3, 5.349936666666667, 16.04981, 0.04148245374559899, 0.04148245374559899, 3, 42.437046, 127.311138, 0.2030878249093402, 0.2119178193804366, 227.03550842042, 2043.31957578378, 0.008424581303096012, 0.008790871140317163)
$newValues[5] = @(3, 5.349936666666667, 16.04981, 0.04148245374559899, 0.04148245374559899, 3, 54.86295166666667, 164.588855, 0.2625535604455031, 0.2739690477507386, 293.5133167630611, 2641.61985086755, 0.01089136592692291, 0.01136490835104581)
$newValues[6] = @(3, 5.349936666666667, 16.04981, 0.04148245374559899, 0.04148245374559899, 2, 26.120138, 52.24027600000001, 0.1250012079717284, 0.08695739860366467, 139.7410840245934, 838.4465041475602, 0.005185356827831223, 0.003607206265414133)
$newValues[7] = @(3, 1.36117, 4.08351, 0.01055426915924182, 0.01055426915924182, 3, 33.67810466666666, 101.034314, 0.1611708087274131, 0.1681783058563055, 45.84162572912667, 412.57463156214, 0.001701040095921799, 0.001774999106752743)
$newValues[8] = @(3, 1.36117, 4.08351, 0.01055426915924182, 0.01055426915924182, 3, 51.86084433333334, 155.582533, 0.2481865979460151, 0.2589774284088547, 70.59142548120336, 635.3228293308301, 0.002619428156438777, 0.002733317485595331)
$newValues[9] = @(3, 1.36117, 4.08351, 0.01055426915924182, 0.01055426915924182, 3, 42.437046, 127.311138, 0.2030878249093402, 0.2119178193804366, 57.76403390382001, 519.87630513438, 0.002143443567058152, 0.002236637705380721)
$newValues[10] = @(3, 1.36117, 4.08351, 0.01055426915924182, 0.01055426915924182, 3, 54.86295166666667, 164.588855, 0.2625535604455031, 0.2739690477507386, 74.67780392011667, 672.10023528105, 0.002771060945659107, 0.00289154307126247)
$newValues[11] = @(3, 1.36117, 4.08351, 0.01055426915924182, 0.01055426915924182, 2, 26.120138, 52.24027600000001, 0.1250012079717284, 0.08695739860366467, 35.55394824146001, 213.3236894487601, 0.001319296394163986, 0.0009177717902505556)
$newValues[12] = @(3, 44.98903266666667, 134.967098, 0.3488369270391816, 0.3488369270391816, 3, 33.67810466666666, 101.034314, 0.1611708087274131, 0.1681783058563055, 1515.145351000086, 13636.30815900077, 0.05622232964489051, 0.0586668034095692)
$newValues[13] = @(3, 44.98903266666667, 134.967098, 0.3488369270391816, 0.3488369270391816, 3, 51.86084433333334, 155.582533, 0.2481865979460151, 0.2589774284088547, 2333.169219833249, 20998.52297849924, 0.08657665015979679, 0.09034089029865451)
$newValues[14] = @(3, 44.98903266666667, 134.967098, 0.3488369270391816, 0.3488369270391816, 3, 42.437046, 127.311138, 0.2030878249093402, 0.2119178193804366, 1909.201648770836, 17182.81483893753, 0.07084453276044561, 0.07392476089751582)
$newValues[15] = @(3, 44.98903266666667, 134.967098, 0.3488369270391816, 0.3488369270391816, 3, 54.86295166666667, 164.588855, 0.2625535604455031, 0.2739690477507386, 2468.231124721422, 22214.08012249279, 0.09158837720900533, 0.09557052072121845)
$newValues[16] = @(3, 44.98903266666667, 134.967098, 0.3488369270391816, 0.3488369270391816, 2, 26.120138, 52.24027600000001, 0.1250012079717284, 0.08695739860366467, 1175.119741739842, 7050.718450439051, 0.04360503726504339, 0.0303339517122236)
$newValues[17] = @(3, 77.26852533333333, 231.805576, 0.5991263500559777, 0.5991263500559777, 3, 33.67810466666666, 101.034314, 0.1611708087274131, 0.1681783058563055, 2602.257483614984, 23420.31735253486, 0.09656167836842515, 0.1007600545462862)
$newValues[18] = @(3, 77.26852533333333, 231.805576, 0.5991263500559777, 0.5991263500559777, 3, 51.86084433333334, 155.582533, 0.2481865979460151, 0.2589774284088547, 4007.210964178224, 36064.89867760401, 0.1486951305602065, 0.1551602014294804)
$newValues[19] = @(3, 77.26852533333333, 231.805576, 0.5991263500559777, 0.5991263500559777, 3, 42.437046, 127.311138, 0.2030878249093402, 0.2119178193804366, 3279.047963922832, 29511.43167530549, 0.1216752672787405, 0.1269655496372229)
$newValues[20] = @(3, 77.26852533333333, 231.805576, 0.5991263500559777, 0.5991263500559777, 3, 54.86295166666667, 164.588855, 0.2625535604455031, 0.2739690477507386, 4239.179370717276, 38152.61433645548, 0.1573027563639158, 0.1641420756072119)
$newValues[21] = @(3, 77.26852533333333, 231.805576, 0.5991263500559777, 0.5991263500559777, 2, 26.120138, 52.24027600000001, 0.1250012079717284, 0.08695739860366467, 2018.264544763163, 12109.58726857898, 0.05209846883577638, 0.05209846883577638)

foreach ($r in $newValues.Keys) {
    $rowVals = $newValues[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowVals[$i]
    }
}
